$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Phase A: delete the three trailing sub-bullet paragraphs that are
# removed entirely by the edit:
#   "I deem any information collect to be damaging or harmful..."
#   "I wish to withdraw from this project."
#   "I no longer wish for group-Q to work on this project."
# ------------------------------------------------------------------
$target = $d.Content.Find
$target.Text = "I no longer wish for group-Q to work on this project."
$target.Execute() | Out-Null
if ($target.Found) {
    $p = $target.Parent.Paragraphs(1)
    $p.Range.Delete()
}

$target = $d.Content.Find
$target.Text = "I wish to withdraw from this project."
$target.Execute() | Out-Null
if ($target.Found) {
    $p = $target.Parent.Paragraphs(1)
    $p.Range.Delete()
}

$target = $d.Content.Find
$target.Text = "I deem any information collect"
$target.Execute() | Out-Null
if ($target.Found) {
    $p = $target.Parent.Paragraphs(1)
    $p.Range.Delete()
}

# ------------------------------------------------------------------
# Phase B: merge "...during the course of this meeting if:" with the
# following "Any of the conditions laid out in..." paragraph by
# deleting the paragraph mark between them. Word gives the merged
# paragraph the *second* paragraph's formatting, so restore the
# original (level-0 bullet) list level afterwards.
# ------------------------------------------------------------------
$target = $d.Content.Find
$target.Text = "Any of the conditions laid out in"
$target.Execute() | Out-Null
$nextPara = $target.Parent.Paragraphs(1)
$joinRange = $d.Range($nextPara.Range.Start - 1, $nextPara.Range.Start)
$joinRange.Delete()

$target = $d.Content.Find
$target.Text = "I (The Participant) retain the right"
$target.Execute() | Out-Null
$mergedPara = $target.Parent.Paragraphs(1)
$mergedPara.Range.ListFormat.ListLevelNumber = 1

# ------------------------------------------------------------------
# Phase C: text edits inside the merged paragraph.
#   " of any data ... this meeting if:Any of the conditions ..."
#     -> " of any data ... this meeting pertaining to myself if any of the conditions ..."
#   each Find/Replace below lands on its own run, approximating the
#   run layout produced by the original incremental Word edits.
# ------------------------------------------------------------------
$d.Content.Find.Execute("meeting if:Any", $true, $false, $false, $false, $false, $true, 1, $false, "meeting~if:Any", 2) | Out-Null
$d.Content.Find.Execute("meeting~if:Any", $true, $false, $false, $false, $false, $true, 1, $false, "meeting", 2) | Out-Null

$find = $d.Content.Find
$find.Text = "meeting"
$find.Execute() | Out-Null
$insertPoint = $d.Range($find.Parent.End, $find.Parent.End)
$insertPoint.InsertAfter(" pertaining to myself")

$find = $d.Content.Find
$find.Text = " pertaining to myself"
$find.Execute() | Out-Null
$afterPertaining = $d.Range($find.Parent.End, $find.Parent.End)
$afterPertaining.InsertAfter(" ")

$find = $d.Content.Find
$find.Text = "if:Any"
$find.Execute() | Out-Null
$find.Parent.Text = "if a"
$find2 = $d.Content.Find
$find2.Text = "if a"
$find2.Execute() | Out-Null
$afterIfA = $d.Range($find2.Parent.End, $find2.Parent.End)
$afterIfA.InsertAfter("ny of the conditions laid out in")

$find = $d.Content.Find
$find.Text = " are violated."
$find.Execute() | Out-Null
$find.Parent.Text = " a"
$find2 = $d.Content.Find
$find2.Text = "(section 1) a"
$find2.Execute() | Out-Null
$endA = $d.Range($find2.Parent.End, $find2.Parent.End)
$endA.InsertAfter("re violated.")

# ------------------------------------------------------------------
# Phase D: move the "_GoBack" bookmark from the end of the last
# paragraph to the end of the merged paragraph (right after
# "...are violated.").
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$find = $d.Content.Find
$find.Text = "re violated."
$find.Execute() | Out-Null
$newBookmarkRange = $d.Range($find.Parent.End, $find.Parent.End)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
